# Update "想去人数" (want-to-go count) figures that changed between
# crawler runs, across the three worksheets that carry this data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value  = 2596
$ws.Range("F10").Value = 2507
$ws.Range("F18").Value = 322
$ws.Range("F32").Value = 1617
$ws.Range("F33").Value = 1019
$ws.Range("F36").Value = 1114
$ws.Range("F37").Value = 2045
$ws.Range("F40").Value = 541
$ws.Range("F44").Value = 1317
$ws.Range("F48").Value = 65

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 67

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value  = 2596
$ws.Range("F8").Value  = 2507
$ws.Range("F15").Value = 322
$ws.Range("F25").Value = 67
$ws.Range("F30").Value = 1617
$ws.Range("F31").Value = 1019
$ws.Range("F34").Value = 2045
$ws.Range("F40").Value = 541
$ws.Range("F44").Value = 1317
$ws.Range("F48").Value = 65

$wb.Save()
